# "Added Cook Park name change"
#
# The park name shown on the map textbox changes from
# "Historic Mims Park" to "Cook Park". Walk every slide (recursing into
# grouped shapes, since the textbox lives inside a picture/textbox group)
# and update the run of text that still says "Historic Mims Park".

$p = $ppt.ActivePresentation
$msoGroup = 6

function Update-ParkName($shape) {
    if ($shape.Type -eq $msoGroup) {
        $items = $shape.GroupItems
        for ($j = 1; $j -le $items.Count; $j++) {
            Update-ParkName $items.Item($j)
        }
        return
    }

    if ($shape.HasTextFrame) {
        $tr = $shape.TextFrame.TextRange
        if ($tr.Text -like "*Historic Mims Park*") {
            $tr.Text = $tr.Text -replace "Historic Mims Park", "Cook Park"
        }
    }
}

for ($i = 1; $i -le $p.Slides.Count; $i++) {
    $slide = $p.Slides.Item($i)
    for ($k = 1; $k -le $slide.Shapes.Count; $k++) {
        Update-ParkName $slide.Shapes.Item($k)
    }
}
